$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update D13: "-" -> "10000 - 40000"
# ---------------------------------------------------------------------------
$ws.Range("D13").Value = "10000 - 40000"

# ---------------------------------------------------------------------------
# 2. New "Core specifications" block in columns F:G
# ---------------------------------------------------------------------------

# F9 / G9 : plain bordered (reuse existing thin-border style) but left empty
$ws.Range("C4").Copy()
$ws.Range("F9:G9").PasteSpecial(-4122)
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = ""

# F4:G8 values - bordered cells (same visual style as normal data cells)
$ws.Range("F4").Value = "Clock Cycle Time"
$ws.Range("G4").Value = "4ns"
$ws.Range("F5").Value = "Memory Operation Time"
$ws.Range("G5").Value = "3 clk cycles"
$ws.Range("F6").Value = "CPI (R,I-TYPE)"
$ws.Range("G6").Value = 1.13
$ws.Range("F7").Value = "Frequency"
$ws.Range("G7").Value = "250MHz"
$ws.Range("F8").Value = "Memory address space"
$ws.Range("G8").Value = "4KB"

$ws.Range("C4").Copy()
$ws.Range("F4:G8").PasteSpecial(-4122)

# F10:G16 - unbordered, centered filler cells
$ws.Range("F10").HorizontalAlignment = -4108
$ws.Range("F10").VerticalAlignment = -4108
$ws.Range("F10").Copy()
$ws.Range("F10:G16").PasteSpecial(-4122)

# F3:G3 - merged header "Core specifications" with grey fill + outer border
$ws.Range("F3").Value = "Core specifications"
$ws.Range("C3").Copy()
$ws.Range("F3:G3").PasteSpecial(-4122)
$ws.Range("F3:G3").Merge()
$ws.Range("F3").Borders.Item(10).LineStyle = 0
$ws.Range("G3").Borders.Item(7).LineStyle = 0

# ---------------------------------------------------------------------------
# 3. Column widths for the new columns
# ---------------------------------------------------------------------------
$ws.Range("F1").EntireColumn.ColumnWidth = 22.0
$ws.Range("G1").EntireColumn.ColumnWidth = 22.5

# ---------------------------------------------------------------------------
# 4. Update the active selection shown when the sheet is opened
# ---------------------------------------------------------------------------
$ws.Range("I15").Select()

Write-Output "done"
